$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: replace the text of a paragraph's Range, preserving the
# formatting of the paragraph's existing run(s). We insert the new text
# right before the end of the existing range (so it inherits formatting
# from the adjacent run) and then delete the old text, being careful to
# leave the paragraph mark itself untouched.
# ---------------------------------------------------------------------
function Set-ParagraphText($para, [string]$newText) {
    $start = $para.Range.Start
    $end = $para.Range.End - 1
    $insPoint = $d.Range($end, $end)
    $insPoint.InsertBefore($newText)
    $oldRange = $d.Range($start, $end)
    $oldRange.Delete()
}

# 1) Heading paragraph: "Other Protected Persons" -> conditional "Other "
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = '{% if trial_court.address.county != "Cook" %}Other {% endif %}Protected Persons'

# 2) Intro paragraph: "Other protected persons (...)" -> conditional P/Other p,
#    and conditional ", in addition to the petitioner..."
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = '{% if trial_court.address.county == "Cook" %}P{% else %}Other p{% endif %}rotected persons (persons to be included in the Stalking No Contact Order){% if trial_court.address.county != "Cook" %}, in addition to the petitioner and the first two protected persons,{% endif %} are:'

# 3) For-loop / if paragraph just before the per-person table.
$p8 = $d.Paragraphs.Item(8)
Set-ParagraphText $p8 '{% for person in others_protected %}{% if person.include_in_addendum == True %}'

# 4) Address cell in the per-person table: wrap with hide_address / safe address logic.
$p10 = $d.Paragraphs.Item(10)
Set-ParagraphText $p10 '{% if hide_address == True and person.use_safe_address == True %}{{users[0].address.on_one_line(bare=True)}}{% else %}{{person.address.on_one_line(bare=True)}}{% endif %}'

# 5) Closing paragraph: "{% endif %}{% endfor %}"
$p12 = $d.Paragraphs.Item(12)
Set-ParagraphText $p12 '{% endif %}{% endfor %}'

# ---------------------------------------------------------------------
# Table width adjustments on the second (per-person) table: shrink the
# address column (and overall table width) to make room.
# ---------------------------------------------------------------------
$t2 = $d.Tables.Item(2)
$t2.PreferredWidth = 458.75
$t2.Columns.Item(2).Width = 234.45
